$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 10
$ws.Range("N3").Value = 10

$ws.Range("N4").Select()
